$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 247.4
$ws.Range("I9").Value = 287.25
$ws.Range("K9").Value = 287.25
$ws.Range("M9").Value = -118.25
# Row 19
$ws.Range("H19").Value = 666.3333
$ws.Range("I19").Value = 499.5
$ws.Range("K19").Value = 499.5
$ws.Range("M19").Value = -324.5
# Row 28
$ws.Range("H28").Value = 900
$ws.Range("I28").Value = 925.15
$ws.Range("K28").Value = 925.15
$ws.Range("M28").Value = -440.15
# Row 33
$ws.Range("H33").Value = 230.6842
$ws.Range("I33").Value = 212.1
$ws.Range("K33").Value = 212.1
$ws.Range("M33").Value = 16.90000000000001
# Row 51
$ws.Range("H51").Value = 87231070
$ws.Range("I51").Value = 149536620
$ws.Range("K51").Value = 149536620
$ws.Range("M51").Value = -149536136
# Row 103
$ws.Range("H103").Value = 399.48
$ws.Range("I103").Value = 237
$ws.Range("J103").Value = 475.94116
$ws.Range("K103").Value = 711
$ws.Range("L103").Value = 1427.82348
$ws.Range("M103").Value = -125
$ws.Range("N103").Value = -2599.82348
# Row 129
$ws.Range("H129").Value = 55556680
$ws.Range("J129").Value = 2401.5
$ws.Range("L129").Value = 7204.5
$ws.Range("N129").Value = -17204.5
# Row 139
$ws.Range("H139").Value = 83199.60000000001
$ws.Range("J139").Value = 83199.60000000001
$ws.Range("L139").Value = 83199.60000000001
$ws.Range("N139").Value = -93479.60000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 17
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10346
# Row 32
$ws.Range("H32").Value = 6931.553
$ws.Range("I32").Value = 7063.9316
$ws.Range("J32").Value = 4990
$ws.Range("K32").Value = 7063.9316
$ws.Range("L32").Value = 4990
$ws.Range("M32").Value = -6776.9316
$ws.Range("N32").Value = -5564
# Row 74
$ws.Range("H74").Value = 2608.2683
$ws.Range("I74").Value = 2139.9429
$ws.Range("J74").Value = 5340.1665
$ws.Range("K74").Value = 2139.9429
$ws.Range("L74").Value = 5340.1665
$ws.Range("M74").Value = -1265.9429
$ws.Range("N74").Value = -7088.1665
# Row 77
$ws.Range("H77").Value = 2608.2683
$ws.Range("I77").Value = 2139.9429
$ws.Range("J77").Value = 5340.1665
$ws.Range("K77").Value = 10699.7145
$ws.Range("L77").Value = 26700.8325
$ws.Range("M77").Value = -6331.7145
$ws.Range("N77").Value = -35436.8325
# Row 88
$ws.Range("H88").Value = 1531.5714
$ws.Range("J88").Value = 1447.7646
$ws.Range("L88").Value = 1447.7646
$ws.Range("N88").Value = -2259.7646
# Row 91
$ws.Range("H91").Value = 1531.5714
$ws.Range("J91").Value = 1447.7646
$ws.Range("L91").Value = 1447.7646
$ws.Range("N91").Value = -4255.7646
# Row 97
$ws.Range("H97").Value = 421.83334
$ws.Range("I97").Value = 421.83334
$ws.Range("K97").Value = 421.83334
$ws.Range("M97").Value = 74.16665999999998
# Row 98
$ws.Range("H98").Value = 76000
$ws.Range("J98").Value = 76000
$ws.Range("L98").Value = 76000
$ws.Range("N98").Value = -81990
# Row 110
$ws.Range("H110").Value = 1646.6666
$ws.Range("I110").Value = 1341.75
$ws.Range("J110").Value = 2866.3333
$ws.Range("K110").Value = 1341.75
$ws.Range("L110").Value = 2866.3333
$ws.Range("M110").Value = 703.25
$ws.Range("N110").Value = -6956.3333
# Row 132
$ws.Range("H132").Value = 2081.0513
$ws.Range("I132").Value = 1850.5
$ws.Range("K132").Value = 5551.5
$ws.Range("M132").Value = -3021.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 37249.25
$ws.Range("I26").Value = 37249.25
$ws.Range("K26").Value = 37249.25
$ws.Range("M26").Value = -36957.25
# Row 35
$ws.Range("H35").Value = 75558
$ws.Range("J35").Value = 75558
$ws.Range("L35").Value = 75558
$ws.Range("N35").Value = -76178
# Row 80
$ws.Range("H80").Value = 790.125
$ws.Range("I80").Value = 551.1667
$ws.Range("J80").Value = 933.5
$ws.Range("K80").Value = 551.1667
$ws.Range("L80").Value = 933.5
$ws.Range("M80").Value = 446.8333
$ws.Range("N80").Value = -2929.5
# Row 83
$ws.Range("H83").Value = 790.125
$ws.Range("I83").Value = 551.1667
$ws.Range("J83").Value = 933.5
$ws.Range("K83").Value = 2755.8335
$ws.Range("L83").Value = 4667.5
$ws.Range("M83").Value = 2236.1665
$ws.Range("N83").Value = -14651.5
# Row 86
$ws.Range("H86").Value = 8883.875
$ws.Range("I86").Value = 11829.909
$ws.Range("J86").Value = 2402.6
$ws.Range("K86").Value = 11829.909
$ws.Range("L86").Value = 2402.6
$ws.Range("M86").Value = -10706.909
$ws.Range("N86").Value = -4648.6
# Row 89
$ws.Range("H89").Value = 8883.875
$ws.Range("I89").Value = 11829.909
$ws.Range("J89").Value = 2402.6
$ws.Range("K89").Value = 59149.545
$ws.Range("L89").Value = 12013
$ws.Range("M89").Value = -53533.545
$ws.Range("N89").Value = -23245
# Row 97
$ws.Range("H97").Value = 3584.375
$ws.Range("I97").Value = 3584.375
$ws.Range("K97").Value = 3584.375
$ws.Range("M97").Value = -2593.375
# Row 99
$ws.Range("H99").Value = 2931.76
$ws.Range("I99").Value = 2164.4119
$ws.Range("K99").Value = 2164.4119
$ws.Range("M99").Value = -666.4119000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2013.7
$ws.Range("I16").Value = 1517.125
$ws.Range("K16").Value = 1517.125
$ws.Range("M16").Value = -1230.125
# Row 22
$ws.Range("H22").Value = 611.7646999999999
$ws.Range("I22").Value = 566.6667
$ws.Range("J22").Value = 720
$ws.Range("K22").Value = 566.6667
$ws.Range("L22").Value = 720
$ws.Range("M22").Value = -216.6667
$ws.Range("N22").Value = -1420
# Row 31
$ws.Range("H31").Value = 2291.5386
$ws.Range("I31").Value = 1405.375
$ws.Range("K31").Value = 1405.375
$ws.Range("M31").Value = -1110.375
# Row 34
$ws.Range("H34").Value = 2291.5386
$ws.Range("I34").Value = 1405.375
$ws.Range("K34").Value = 1405.375
$ws.Range("M34").Value = -1203.375
# Row 113
$ws.Range("H113").Value = 2013.7
$ws.Range("I113").Value = 1517.125
$ws.Range("K113").Value = 1517.125
$ws.Range("M113").Value = 652.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 96.71429000000001
$ws.Range("I40").Value = 39.272728
$ws.Range("J40").Value = 307.33334
$ws.Range("K40").Value = 157.090912
$ws.Range("L40").Value = 1229.33336
$ws.Range("M40").Value = -88.090912
$ws.Range("N40").Value = -1367.33336
# Row 107
$ws.Range("H107").Value = 7993.4287
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 9259
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 27777
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -31617
# Row 117
$ws.Range("H117").Value = 1036.4166
$ws.Range("I117").Value = 697.8570999999999
$ws.Range("J117").Value = 1510.4
$ws.Range("K117").Value = 2093.5713
$ws.Range("L117").Value = 4531.200000000001
$ws.Range("M117").Value = 1348.4287
$ws.Range("N117").Value = -11415.2
# Row 131
$ws.Range("H131").Value = 1730.7368
$ws.Range("J131").Value = 1959
$ws.Range("L131").Value = 5877
$ws.Range("N131").Value = -15957
# Row 137
$ws.Range("H137").Value = 3100.6924
$ws.Range("J137").Value = 3192.8572
$ws.Range("L137").Value = 9578.571599999999
$ws.Range("N137").Value = -19778.5716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 128
$ws.Range("H128").Value = 89999.25
$ws.Range("J128").Value = 89999.25
$ws.Range("L128").Value = 89999.25
$ws.Range("N128").Value = -99959.25
# Row 132
$ws.Range("H132").Value = 11909464
$ws.Range("I132").Value = 18523044
$ws.Range("K132").Value = 55569132
$ws.Range("M132").Value = -55566602

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 4284.6665
$ws.Range("I4").Value = 4154.5
$ws.Range("J4").Value = 4545
$ws.Range("K4").Value = 4154.5
$ws.Range("L4").Value = 4545
$ws.Range("M4").Value = -4041.5
$ws.Range("N4").Value = -4771
# Row 22
$ws.Range("H22").Value = 2727.15
$ws.Range("I22").Value = 1361.25
$ws.Range("J22").Value = 3068.625
$ws.Range("K22").Value = 1361.25
$ws.Range("L22").Value = 3068.625
$ws.Range("M22").Value = -1066.25
$ws.Range("N22").Value = -3658.625
# Row 27
$ws.Range("H27").Value = 2727.15
$ws.Range("I27").Value = 1361.25
$ws.Range("J27").Value = 3068.625
$ws.Range("K27").Value = 1361.25
$ws.Range("L27").Value = 3068.625
$ws.Range("M27").Value = -1254.25
$ws.Range("N27").Value = -3282.625
# Row 28
$ws.Range("H28").Value = 4284.6665
$ws.Range("I28").Value = 4154.5
$ws.Range("J28").Value = 4545
$ws.Range("K28").Value = 4154.5
$ws.Range("L28").Value = 4545
$ws.Range("M28").Value = -3922.5
$ws.Range("N28").Value = -5009
# Row 37
$ws.Range("H37").Value = 4284.6665
$ws.Range("I37").Value = 4154.5
$ws.Range("J37").Value = 4545
$ws.Range("K37").Value = 4154.5
$ws.Range("L37").Value = 4545
$ws.Range("M37").Value = -4047.5
$ws.Range("N37").Value = -4759
# Row 97
$ws.Range("H97").Value = 56637.6
$ws.Range("J97").Value = 56637.6
$ws.Range("L97").Value = 56637.6
$ws.Range("N97").Value = -58619.6
# Row 98
$ws.Range("H98").Value = 72333
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
# Row 105
$ws.Range("H105").Value = 43300
$ws.Range("J105").Value = 43300
$ws.Range("L105").Value = 43300
$ws.Range("N105").Value = -50288
# Row 136
$ws.Range("H136").Value = 4825.6665
$ws.Range("I136").Value = 4706.467
$ws.Range("J136").Value = 5123.6665
$ws.Range("K136").Value = 14119.401
$ws.Range("L136").Value = 15370.9995
$ws.Range("M136").Value = -11569.401
$ws.Range("N136").Value = -20470.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 54442.75
$ws.Range("J94").Value = 42590.332
$ws.Range("L94").Value = 42590.332
$ws.Range("N94").Value = -44392.332
# Row 114
$ws.Range("H114").Value = 89999.5
$ws.Range("J114").Value = 89999.5
$ws.Range("L114").Value = 89999.5
$ws.Range("N114").Value = -98677.5
# Row 122
$ws.Range("H122").Value = 9837.125
$ws.Range("I122").Value = 9813.857
$ws.Range("K122").Value = 29441.571
$ws.Range("M122").Value = -26991.571

